$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 505.7143
$ws.Range("I118").Value = 423.33334
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 1270.00002
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = 386.9999800000001
$ws.Range("N118").Value = -6314
$ws.Range("H129").Value = 2257.5483
$ws.Range("I129").Value = 764.6667
$ws.Range("J129").Value = 2417.5
$ws.Range("K129").Value = 2294.0001
$ws.Range("L129").Value = 7252.5
$ws.Range("M129").Value = 2705.9999
$ws.Range("N129").Value = -17252.5
$ws.Range("H132").Value = 2719336
$ws.Range("I132").Value = 3473495.8
$ws.Range("J132").Value = 4360.4
$ws.Range("K132").Value = 10420487.4
$ws.Range("L132").Value = 13081.2
$ws.Range("M132").Value = -10417957.4
$ws.Range("N132").Value = -18141.2
$ws.Range("H137").Value = 883
$ws.Range("I137").Value = 758.625
$ws.Range("K137").Value = 2275.875
$ws.Range("M137").Value = 274.125
$ws.Range("H138").Value = 2600.6973
$ws.Range("I138").Value = 877.64105
$ws.Range("J138").Value = 4416.892
$ws.Range("K138").Value = 2632.92315
$ws.Range("L138").Value = 13250.676
$ws.Range("M138").Value = 2507.07685
$ws.Range("N138").Value = -23530.676
$ws.Range("H141").Value = 4903.303
$ws.Range("I141").Value = 3713.2778
$ws.Range("J141").Value = 6331.3335
$ws.Range("K141").Value = 11139.8334
$ws.Range("L141").Value = 18994.0005
$ws.Range("M141").Value = -5959.8334
$ws.Range("N141").Value = -29354.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1906.7142
$ws.Range("I61").Value = 1509.3334
$ws.Range("J61").Value = 2622
$ws.Range("K61").Value = 1509.3334
$ws.Range("L61").Value = 2622
$ws.Range("M61").Value = -1297.3334
$ws.Range("N61").Value = -3046
$ws.Range("H74").Value = 834.73334
$ws.Range("I74").Value = 725.6667
$ws.Range("J74").Value = 1271
$ws.Range("K74").Value = 725.6667
$ws.Range("L74").Value = 1271
$ws.Range("M74").Value = 148.3333
$ws.Range("N74").Value = -3019
$ws.Range("H77").Value = 834.73334
$ws.Range("I77").Value = 725.6667
$ws.Range("J77").Value = 1271
$ws.Range("K77").Value = 3628.3335
$ws.Range("L77").Value = 6355
$ws.Range("M77").Value = 739.6665000000003
$ws.Range("N77").Value = -15091
$ws.Range("H132").Value = 2218.875
$ws.Range("I132").Value = 1318.7142
$ws.Range("K132").Value = 3956.1426
$ws.Range("M132").Value = -1426.1426
$ws.Range("H136").Value = 1906.7142
$ws.Range("I136").Value = 1509.3334
$ws.Range("J136").Value = 2622
$ws.Range("K136").Value = 4528.0002
$ws.Range("L136").Value = 7866
$ws.Range("M136").Value = -1978.0002
$ws.Range("N136").Value = -12966

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 272.88
$ws.Range("I22").Value = 272.88
$ws.Range("K22").Value = 272.88
$ws.Range("M22").Value = -99.88
$ws.Range("H134").Value = 2148.25
$ws.Range("I134").Value = 1888.2858
$ws.Range("J134").Value = 3058.125
$ws.Range("K134").Value = 5664.857400000001
$ws.Range("L134").Value = 9174.375
$ws.Range("M134").Value = -3129.857400000001
$ws.Range("N134").Value = -14244.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1754.8636
$ws.Range("I31").Value = 1352.7333
$ws.Range("J31").Value = 2616.5715
$ws.Range("K31").Value = 1352.7333
$ws.Range("L31").Value = 2616.5715
$ws.Range("M31").Value = -1057.7333
$ws.Range("N31").Value = -3206.5715
$ws.Range("H33").Value = 30333.334
$ws.Range("I33").Value = 1000
$ws.Range("J33").Value = 45000
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 45000
$ws.Range("M33").Value = -621
$ws.Range("N33").Value = -45758
$ws.Range("H34").Value = 1754.8636
$ws.Range("I34").Value = 1352.7333
$ws.Range("J34").Value = 2616.5715
$ws.Range("K34").Value = 1352.7333
$ws.Range("L34").Value = 2616.5715
$ws.Range("M34").Value = -1150.7333
$ws.Range("N34").Value = -3020.5715
$ws.Range("H58").Value = 868.8372000000001
$ws.Range("I58").Value = 630.625
$ws.Range("J58").Value = 1561.8182
$ws.Range("K58").Value = 630.625
$ws.Range("L58").Value = 1561.8182
$ws.Range("M58").Value = -427.625
$ws.Range("N58").Value = -1967.8182
$ws.Range("H132").Value = 8235.736999999999
$ws.Range("I132").Value = 9320.071
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 27960.213
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -25430.213
$ws.Range("N132").Value = -20658.8
$ws.Range("H134").Value = 1355.8846
$ws.Range("I134").Value = 1218.875
$ws.Range("K134").Value = 3656.625
$ws.Range("M134").Value = -1121.625
$ws.Range("H136").Value = 868.8372000000001
$ws.Range("I136").Value = 630.625
$ws.Range("J136").Value = 1561.8182
$ws.Range("K136").Value = 1891.875
$ws.Range("L136").Value = 4685.4546
$ws.Range("M136").Value = 658.125
$ws.Range("N136").Value = -9785.454600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 4296.4287
$ws.Range("J126").Value = 5321.5
$ws.Range("L126").Value = 15964.5
$ws.Range("N126").Value = -25844.5
$ws.Range("H130").Value = 3656.5557
$ws.Range("I130").Value = 1515
$ws.Range("J130").Value = 4268.4287
$ws.Range("K130").Value = 4545
$ws.Range("L130").Value = 12805.2861
$ws.Range("M130").Value = 475
$ws.Range("N130").Value = -22845.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4475.0884
$ws.Range("I132").Value = 4629.483
$ws.Range("J132").Value = 3579.6
$ws.Range("K132").Value = 13888.449
$ws.Range("L132").Value = 10738.8
$ws.Range("M132").Value = -11358.449
$ws.Range("N132").Value = -15798.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2031.9592
$ws.Range("I132").Value = 1298.5428
$ws.Range("J132").Value = 3865.5
$ws.Range("K132").Value = 3895.6284
$ws.Range("L132").Value = 11596.5
$ws.Range("M132").Value = -1365.6284
$ws.Range("N132").Value = -16656.5
$ws.Range("H136").Value = 2040
$ws.Range("I136").Value = 1134
$ws.Range("J136").Value = 2704.4
$ws.Range("K136").Value = 3402
$ws.Range("L136").Value = 8113.200000000001
$ws.Range("M136").Value = -852
$ws.Range("N136").Value = -13213.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1637.25
$ws.Range("I132").Value = 918.4737
$ws.Range("J132").Value = 2287.5715
$ws.Range("K132").Value = 2755.4211
$ws.Range("L132").Value = 6862.7145
$ws.Range("M132").Value = -225.4211
$ws.Range("N132").Value = -11922.7145
$ws.Range("H136").Value = 3906.25
$ws.Range("I136").Value = 993.2273
$ws.Range("K136").Value = 2979.6819
$ws.Range("M136").Value = -429.6819
